# ---------------------------------------------------------------------------
# Applies the "Add analysis to xlsx" commit to exp7/data/readings.xlsx
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix the coupling-constant table (rows 31-34): the denominator used to
#    divide by $B$42 (wrong / unrelated cell) and must instead divide by
#    $J$32 (the "n" value for this particular pendulum pair).
# ---------------------------------------------------------------------------
$ws.Range("L32").Formula = '=K32/$J$32'
$ws.Range("L33").Formula = '=K33/$J$32'
$ws.Range("L34").Formula = '=K34/$J$32'

# L34 loses its bottom border (it now matches the styling of L32/L33).
$ws.Range("L34").Style = $ws.Range("L32").Style

# ---------------------------------------------------------------------------
# 2) Extend the "coupling moment" table (row 36 headers / row 37 data) with
#    two new uncertainty columns: Sphi [deg] and sk.
# ---------------------------------------------------------------------------
$ws.Range("O36").Value = "Sphi [" + [char]0x00B0 + "]"
$ws.Range("P36").Value = "sk"

$ws.Range("O37").Value = 0.2
$ws.Range("P37").Formula = '=SQRT((M37^2*O37^2+L37^2*O37^2)/M37^4)'

# ---------------------------------------------------------------------------
# 3) Add the "exp" dynamic-k / dynamic-tau analysis block next to the
#    existing k-table (rows 40-42).
# ---------------------------------------------------------------------------
$ws.Range("K40").Value = "exp"

$ws.Range("K41").Value = "kdyn"
$ws.Range("L41").Value = "s kdyn"
$ws.Range("M41").Value = "tau"
$ws.Range("N41").Value = "s tau"
$ws.Range("O41").Value = "ts"
$ws.Range("P41").Value = "s ts"

$ws.Range("K42").Formula = '=(O19^2-O23^2)/(O19^2+O23^2)'
$ws.Range("L42").Formula = '=SQRT(16*(P23^2*O19^4*O23^2+P19^2*O19^2*O23^4)/(O19^2+O23^2)^4)'
$ws.Range("M42").Formula = '=1/(1/2*(1/O19+1/O23))'
$ws.Range("O42").Formula = '=1/(1/O23-1/O19)'

# ---------------------------------------------------------------------------
# 4) Row 46 used to contain a duplicate (mis-addressed) copy of row 45's
#    data; it must become an empty row.
# ---------------------------------------------------------------------------
$ws.Range("A46:F46").Clear()

# ---------------------------------------------------------------------------
# 5) Mirror the same two additions (Sphi[deg]/sk header, uncertainty value
#    + formula) on the second "coupling moment" table (rows 47-48).
# ---------------------------------------------------------------------------
$ws.Range("G47").Value = "Sphi [" + [char]0x00B0 + "]"
$ws.Range("H47").Value = "sk"

$ws.Range("G48").Value = 0.2
$ws.Range("H48").Formula = '=SQRT((E48^2*G48^2+D48^2*G48^2)/E48^4)'

# ---------------------------------------------------------------------------
# 6) Mirror the "exp" dynamic-k / dynamic-tau analysis block for the second
#    table (rows 50-52), including the brand new row 52.
# ---------------------------------------------------------------------------
$ws.Range("C50").Value = "exp"

$ws.Range("C51").Value = "kdyn"
$ws.Range("D51").Value = "s kdyn"
$ws.Range("E51").Value = "tau"
$ws.Range("F51").Value = "s tau"
$ws.Range("G51").Value = "ts"
$ws.Range("H51").Value = "s ts"

$ws.Range("C52").Formula = '=(G19^2-G25^2)/(G19^2+G25^2)'
$ws.Range("D52").Formula = '=SQRT(16*(H25^2*G19^4*G25^2+H19^2*G19^2*G25^4)/(G19^2+G25^2)^4)'
$ws.Range("E52").Formula = '=1/(1/2*(1/G19+1/G25))'
$ws.Range("F52").Formula = '=SQRT(4*(H25^2*G19^4+H19^2*G25^4)/(G19+G25)^4)'
$ws.Range("G52").Formula = '=1/(1/G25-1/G19)'

# ---------------------------------------------------------------------------
# 7) Cosmetic: widen column D slightly and update the window / selection
#    state to match how the author last left the sheet.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 8.3

$ws.Range("F52").Select()

$aw = $excel.ActiveWindow
$aw.Zoom = 110
$aw.ScrollRow = 20
$aw.ScrollColumn = 1
$wb.Windows.Item(1).TabRatio = 0.23

# ---------------------------------------------------------------------------
# 8) Cosmetic: consolidate the redundant "GENERAL" number format (165) back
#    onto the default style (164) for every cell that still referenced it.
# ---------------------------------------------------------------------------
$dedupCells = @('G10','H10','G19','H19','O19','P19','O23','P23','G25','H25','O27','P27','O31','P31','G33','H33','G41','H41')
foreach ($c in $dedupCells) {
    $ws.Range($c).Style = "Normal"
}

Write-Host "Edit applied successfully"
